# Insert a new row at position 259, shifting existing rows 259:374 down to 260:375,
# then populate the new row 259 with the new data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 259 (existing data shifts down)
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new record's values
$ws.Range("A259").Value = 10
$ws.Range("B259").Value = "Vega Modelo de Temuco"
$ws.Range("C259").Value = "La Araucanía"
$ws.Range("D259").Value = 45141
$ws.Range("E259").Value = 9
$ws.Range("F259").Value = 100112043
$ws.Range("G259").Value = "Pepino dulce"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 400
$ws.Range("K259").Value = 23000
$ws.Range("L259").Value = 23000
$ws.Range("M259").Value = 23000
$ws.Range("N259").Value = "$/bandeja 18 kilos"
$ws.Range("O259").Value = "Provincia de Limarí"
$ws.Range("P259").Value = 1278
$ws.Range("Q259").Value = 18
$ws.Range("R259").Value = "Hortaliza"
